# "Forgot Password dan update login"
#
# Updates the User-Login sheet (remember-me row) and the
# User - Forgot Password sheet (new npkNotFound failure row / renumbered
# OTP-failure rows), fixes a couple of shared-string typos, and switches
# the active tab / selections to match the new editing focus.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)  # User-Login
$ws2 = $wb.Worksheets.Item(2)  # User - Forgot Password

# ------------------------------------------------------------------
# Sheet 1: User-Login
# ------------------------------------------------------------------

$ws1.Range("A5").Value = 14426
$ws1.Range("B5").Value = "Password"

$ws1.Range("B6").Value = "Password2"
$ws1.Range("C6").Value = "pass"
$ws1.Range("D6").Value = "remember"

# ------------------------------------------------------------------
# Sheet 2: User - Forgot Password
# ------------------------------------------------------------------

# Typo fix
$ws2.Range("K5").Value = "emptyPasswordConfirm"

# New npk-not-found failure case (row 3), right aligned like a "code" cell
$ws2.Range("A3").Value = "1a2b3c"
$ws2.Range("A3").HorizontalAlignment = -4152   # xlRight

$ws2.Range("K10").Value = "npkNotFound"

$ws2.Range("J2").Value = "fail1"
$ws2.Range("J3").Value = "fail1"
$ws2.Range("J10").Value = "fail1"

# Row 10 loses its OTP-digits sample data ...
$ws2.Range("A10").Value = 55555
$ws2.Range("B10").ClearContents()
$ws2.Range("C10").ClearContents()
$ws2.Range("D10").ClearContents()
$ws2.Range("E10").ClearContents()
$ws2.Range("F10").ClearContents()
$ws2.Range("G10").ClearContents()
$ws2.Range("H10").ClearContents()
$ws2.Range("I10").ClearContents()

# ... which reappears (unchanged) on row 11, now tagged as the wrongOTP case
$ws2.Range("D11").Value = 1
$ws2.Range("E11").Value = 9
$ws2.Range("F11").Value = 9
$ws2.Range("G11").Value = 9
$ws2.Range("H11").Value = 8
$ws2.Range("I11").Value = 9
$ws2.Range("J11").Value = "fail"
$ws2.Range("K11").Value = "wrongOTP"

# New row 12: a second OTP-mismatch sample that now passes
$ws2.Range("A12").Value = 14426
$ws2.Range("B12").Value = "Password2"
$ws2.Range("C12").Value = "Password2"
$ws2.Range("D12").Value = 2
$ws2.Range("E12").Value = 4
$ws2.Range("F12").Value = 9
$ws2.Range("G12").Value = 9
$ws2.Range("H12").Value = 2
$ws2.Range("I12").Value = 1
$ws2.Range("J12").Value = "pass"

# ------------------------------------------------------------------
# Selections / active tab
# ------------------------------------------------------------------

$ws1.Activate()
[void]$ws1.Rows("6").Select()

$ws2.Activate()
[void]$ws2.Range("A16:A17").Select()
